$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data update: add 3 to every numeric cell in columns B:BX across rows 2-5.
# Rows 2 and 3 are already literal numbers; rows 4 and 5 are formulas whose
# cached results need to be "flattened" into literal values. Reading the
# range through .Value2, adjusting it in memory, and writing it back as
# .Value2 achieves both: it bumps every number by 3 AND replaces any
# formula cells with plain literal values (since we are writing values,
# not formulas).
$rng = $ws.Range("B2:BX5")
$vals = $rng.Value2
for ($r = 1; $r -le $vals.GetLength(0); $r++) {
  for ($c = 1; $c -le $vals.GetLength(1); $c++) {
    $vals[$r, $c] = $vals[$r, $c] + 3
  }
}
$rng.Value2 = $vals

# --- Column BY (boolean helper column) keeps its TRUE/FALSE value but the
# formulas (=FALSE / =TRUE) are removed, leaving plain boolean literals.
$ws.Range("BY2").Value2 = $false
$ws.Range("BY3").Value2 = $false
$ws.Range("BY4").Value2 = $false
$ws.Range("BY5").Value2 = $true

# --- Selection moves from BW2 to AG24.
$ws.Range("AG24").Select() | Out-Null
